# Apply the "Add files via upload" update to TC_Template.xlsx:
#  - Fill in the previously blank "Automated" (column G) cells on the
#    "Registration" sheet with Y/N values.
#  - Make "Registration" the active/selected sheet (it previously was "1.2"),
#    scroll it so column C is at the left edge, and leave the selection on I10.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Registration")

# --- Column G ("Automated") values for rows 2-11 on the Registration sheet ---
$ws1.Range("G2").Value  = "N"
$ws1.Range("G3").Value  = "N"
$ws1.Range("G4").Value  = "N"
$ws1.Range("G5").Value  = "N"
$ws1.Range("G6").Value  = "N"
$ws1.Range("G7").Value  = "Y"
$ws1.Range("G8").Value  = "Y"
$ws1.Range("G9").Value  = "N"
$ws1.Range("G10").Value = "N"
$ws1.Range("G11").Value = "N"

# --- Make "Registration" the active sheet/tab (was "1.2") ---
# Activating it also clears the tabSelected flag that was previously on "1.2".
$ws1.Activate()

# Scroll the view so column C is the left-most visible column (topLeftCell=C5)
$win = $excel.ActiveWindow
$win.ScrollColumn = 3
$win.ScrollRow = 5

# Update the selection on the Registration sheet to I10
$ws1.Range("I10").Select()
